$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-16 Friday" "2026-01-17 Saturday"

Replace-Text "94×79=7426" "20×65=1300"
Replace-Text "57×23=1311" "99×85=8415"
Replace-Text "88×62=5456" "59×90=5310"
Replace-Text "97×54=5238" "15×51=765"
Replace-Text "51×84=4284" "89×87=7743"

Replace-Text "67×81=5427" "96×58=5568"
Replace-Text "57×48=2736" "22×91=2002"
Replace-Text "73×76=5548" "42×97=4074"
Replace-Text "95×86=8170" "25×49=1225"
Replace-Text "28×40=1120" "27×74=1998"

Replace-Text "27×59=1593" "29×59=1711"
Replace-Text "19×58=1102" "94×19=1786"
Replace-Text "93×31=2883" "21×97=2037"
Replace-Text "30×18=540" "65×35=2275"
Replace-Text "57×60=3420" "63×95=5985"

Replace-Text "60×81=4860" "73×93=6789"
Replace-Text "43×92=3956" "30×56=1680"
Replace-Text "64×66=4224" "66×12=792"
Replace-Text "99×64=6336" "87×91=7917"
Replace-Text "40×19=760" "94×89=8366"

Replace-Text "95×87=8265" "55×33=1815"
Replace-Text "39×21=819" "33×27=891"
Replace-Text "84×55=4620" "51×87=4437"
Replace-Text "83×37=3071" "29×40=1160"
Replace-Text "84×35=2940" "71×85=6035"
